$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 8.841467
$ws.Range("H2").Value = 26.524401
$ws.Range("I2").Value = 0.5917001192060068
$ws.Range("J2").Value = 0.5917001192060067
$ws.Range("O2").Value = 0.921725411846598
$ws.Range("P2").Value = 0.9217254118465981
$ws.Range("Q2").Value = 56.11019236746233
$ws.Range("R2").Value = 504.991731307161
$ws.Range("S2").Value = 0.5453850360648378
$ws.Range("T2").Value = 0.5453850360648377

# Row 3
$ws.Range("G3").Value = 8.841467
$ws.Range("H3").Value = 26.524401
$ws.Range("I3").Value = 0.5917001192060068
$ws.Range("J3").Value = 0.5917001192060067
$ws.Range("M3").Value = 0.5389353333333333
$ws.Range("N3").Value = 1.616806
$ws.Range("O3").Value = 0.07827458815340194
$ws.Range("P3").Value = 0.07827458815340194
$ws.Range("Q3").Value = 4.764978964800666
$ws.Range("R3").Value = 42.884810683206
$ws.Range("S3").Value = 0.04631508314116901
$ws.Range("T3").Value = 0.04631508314116901

# Row 4
$ws.Range("I4").Value = 0.2746155987184545
$ws.Range("J4").Value = 0.2746155987184545
$ws.Range("O4").Value = 0.921725411846598
$ws.Range("P4").Value = 0.9217254118465981
$ws.Range("S4").Value = 0.2531201758282676
$ws.Range("T4").Value = 0.2531201758282676

# Row 5
$ws.Range("I5").Value = 0.2746155987184545
$ws.Range("J5").Value = 0.2746155987184545
$ws.Range("M5").Value = 0.5389353333333333
$ws.Range("N5").Value = 1.616806
$ws.Range("O5").Value = 0.07827458815340194
$ws.Range("P5").Value = 0.07827458815340194
$ws.Range("Q5").Value = 2.211487726342666
$ws.Range("R5").Value = 19.903389537084
$ws.Range("S5").Value = 0.02149542289018692
$ws.Range("T5").Value = 0.02149542289018692

# Row 6
$ws.Range("G6").Value = 1.997574666666667
$ws.Range("H6").Value = 5.992724
$ws.Range("I6").Value = 0.1336842820755386
$ws.Range("J6").Value = 0.1336842820755386
$ws.Range("O6").Value = 0.921725411846598
$ws.Range("P6").Value = 0.9217254118465981
$ws.Range("Q6").Value = 12.67711555277378
$ws.Range("R6").Value = 114.094039974964
$ws.Range("S6").Value = 0.1232201999534926
$ws.Range("T6").Value = 0.1232201999534926

# Row 7
$ws.Range("G7").Value = 1.997574666666667
$ws.Range("H7").Value = 5.992724
$ws.Range("I7").Value = 0.1336842820755386
$ws.Range("J7").Value = 0.1336842820755386
$ws.Range("M7").Value = 0.5389353333333333
$ws.Range("N7").Value = 1.616806
$ws.Range("O7").Value = 0.07827458815340194
$ws.Range("P7").Value = 0.07827458815340194
$ws.Range("Q7").Value = 1.076563568838222
$ws.Range("R7").Value = 9.689072119543999
$ws.Range("S7").Value = 0.010464082122046
$ws.Range("T7").Value = 0.010464082122046
